# "fixed & subbitted 2.7 working on project fixing animation product diary"
#
# A second "LinearSearchMethod_fixed" test case (row 4, columns M:Q) is
# added, mirroring the first one already present in row 3: a test number,
# two image hyperlinks ("test" / "outcome") with the usual red/green
# pass-fail fill, and a boolean result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New test case #2 for the "LinearSearchMethod_fixed" block ---------

# Test number (plain, unstyled number like M3).
$ws.Range("M4").Value = 2

# Add the two external image hyperlinks first (while the cells are still
# plain/unstyled), then stamp the same look-and-feel as row 3's
# equivalent cells (N3:Q3 -> fill colour + wrap-text + hyperlink font)
# on top via a format-only paste, and finally (re)write the cell values/
# text so the paste doesn't clobber them.
$ws.Hyperlinks.Add($ws.Range("N4"), "LinearSearchMethod_Fixed test 2.png")
$ws.Hyperlinks.Add($ws.Range("P4"), "LinearSearchMethod_Fixed test outcome 2.png")

$ws.Range("N3:Q3").Copy()
$ws.Range("N4:Q4").PasteSpecial(-4122)

$ws.Range("N4").Value = "Testing if this works (Click this for image)"
$ws.Range("O4").Value = "Should print 'watermelon 1' on the first line then print none on the second line"
$ws.Range("P4").Value = "Here is the outcome (don't worry about the falses)"
$ws.Range("Q4").Value = $true

# --- Selection left where the author ended up after the edit -----------
$ws.Range("N4").Select()
